$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "Groups" columns
$ws.Range("F1").Value = "СВОБОДНЫХ МЕСТ"
$ws.Range("G1").Value = "МУЖЧИН"
$ws.Range("H1").Value = "ЖЕНЩИН"

# Set column widths to match the target layout (values chosen so the
# engine's char-width/pixel quantization lands on the closest cell to the
# authored widths of 16.75 / 8.26 / 8.82)
$ws.Columns.Item(6).ColumnWidth = 15.75
$ws.Columns.Item(7).ColumnWidth = 7.42
$ws.Columns.Item(8).ColumnWidth = 8

# Update the selected cell to D1
$ws.Range("D1").Select()
